$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.285.12'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.738.13'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.56'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.81'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").Value = '3.735.53'
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  -3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '4.370.26'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '3.748.35'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '68.263.79'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.73'
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.111'
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.61'
$ws.Range("E21").Value = '  +1.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.01'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.691'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.86'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000146'
$ws.Range("E25").Value = '  +5.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.84'
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.02'
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '3.889.89'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.26'
$ws.Range("E32").Value = '  -1.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.73'
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("E34").Value = '  -2.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.12'
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D37").Value = '3.697.46'
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.100'
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.41'
$ws.Range("E39").Value = '  -4.43%  '
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.995'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.76'
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.300'
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.03'
$ws.Range("E46").Value = '  +10.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.68'
$ws.Range("E47").Value = '  +3.98%  '
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.43'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '387.57'
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.28'
$ws.Range("E51").Value = '  +0.18%  '
